$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 12502125
$ws.Range("J19").Value = 12502125
$ws.Range("L19").Value = 12502125
$ws.Range("N19").Value = -12502475
$ws.Range("H33").Value = 258.66666
$ws.Range("I33").Value = 235.05556
$ws.Range("J33").Value = 400.33334
$ws.Range("K33").Value = 235.05556
$ws.Range("L33").Value = 400.33334
$ws.Range("M33").Value = -6.055560000000014
$ws.Range("N33").Value = -858.33334
$ws.Range("H51").Value = 6710.625
$ws.Range("I51").Value = 9104.1
$ws.Range("J51").Value = 5622.6816
$ws.Range("K51").Value = 9104.1
$ws.Range("L51").Value = 5622.6816
$ws.Range("M51").Value = -8620.1
$ws.Range("N51").Value = -6590.6816
$ws.Range("H88").Value = 4831.8
$ws.Range("J88").Value = 6022.75
$ws.Range("L88").Value = 6022.75
$ws.Range("N88").Value = -6834.75
$ws.Range("H91").Value = 4831.8
$ws.Range("J91").Value = 6022.75
$ws.Range("L91").Value = 6022.75
$ws.Range("N91").Value = -8830.75
$ws.Range("H98").Value = 1792.28
$ws.Range("I98").Value = 1661.1765
$ws.Range("J98").Value = 2070.875
$ws.Range("K98").Value = 1661.1765
$ws.Range("L98").Value = 2070.875
$ws.Range("M98").Value = -163.1765
$ws.Range("N98").Value = -5066.875
$ws.Range("H122").Value = 1792.28
$ws.Range("I122").Value = 1661.1765
$ws.Range("J122").Value = 2070.875
$ws.Range("K122").Value = 4983.529500000001
$ws.Range("L122").Value = 6212.625
$ws.Range("M122").Value = -2533.529500000001
$ws.Range("N122").Value = -11112.625
$ws.Range("H129").Value = 2110.7058
$ws.Range("I129").Value = 833.6667
$ws.Range("K129").Value = 2501.0001
$ws.Range("M129").Value = 2498.9999
$ws.Range("H137").Value = 14288376
$ws.Range("I137").Value = 22729552
$ws.Range("K137").Value = 68188656
$ws.Range("M137").Value = -68186106
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 118334330
$ws.Range("I61").Value = 175001000
$ws.Range("K61").Value = 175001000
$ws.Range("M61").Value = -175000788
$ws.Range("H132").Value = 2044213.5
$ws.Range("I132").Value = 3403.7856
$ws.Range("K132").Value = 10211.3568
$ws.Range("M132").Value = -7681.356800000001
$ws.Range("H136").Value = 118334330
$ws.Range("I136").Value = 175001000
$ws.Range("K136").Value = 525003000
$ws.Range("M136").Value = -525000450
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H26").Value = 34048.125
$ws.Range("I26").Value = 24626.572
$ws.Range("J26").Value = 99999
$ws.Range("K26").Value = 24626.572
$ws.Range("L26").Value = 99999
$ws.Range("M26").Value = -24334.572
$ws.Range("N26").Value = -100583
$ws.Range("H94").Value = 3078.476
$ws.Range("I94").Value = 3116.6365
$ws.Range("K94").Value = 3116.6365
$ws.Range("M94").Value = -2665.6365
$ws.Range("H99").Value = 3031.8
$ws.Range("I99").Value = 3039.875
$ws.Range("K99").Value = 3039.875
$ws.Range("M99").Value = -1541.875
$ws.Range("H134").Value = 5885639
$ws.Range("I134").Value = 3186.8
$ws.Range("K134").Value = 9560.400000000001
$ws.Range("M134").Value = -7025.400000000001
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 35717908
$ws.Range("J31").Value = 4090.375
$ws.Range("L31").Value = 4090.375
$ws.Range("N31").Value = -4680.375
$ws.Range("H34").Value = 35717908
$ws.Range("J34").Value = 4090.375
$ws.Range("L34").Value = 4090.375
$ws.Range("N34").Value = -4494.375
$ws.Range("H86").Value = 13311.091
$ws.Range("I86").Value = 13311.091
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 13311.091
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -12188.091
$ws.Range("N86").Value = $null
$ws.Range("H89").Value = 13311.091
$ws.Range("I89").Value = 13311.091
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 66555.455
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -60939.455
$ws.Range("N89").Value = $null
$ws.Range("H99").Value = 66999.39999999999
$ws.Range("I99").Value = 14999
$ws.Range("J99").Value = 145000
$ws.Range("K99").Value = 14999
$ws.Range("L99").Value = 145000
$ws.Range("M99").Value = -13501
$ws.Range("N99").Value = -147996
$ws.Range("H122").Value = 3806.818
$ws.Range("I122").Value = 3548.125
$ws.Range("K122").Value = 10644.375
$ws.Range("M122").Value = -8194.375
$ws.Range("H126").Value = 66999.39999999999
$ws.Range("I126").Value = 14999
$ws.Range("J126").Value = 145000
$ws.Range("K126").Value = 44997
$ws.Range("L126").Value = 435000
$ws.Range("M126").Value = -42527
$ws.Range("N126").Value = -439940
$ws.Range("H132").Value = 2552.375
$ws.Range("I132").Value = 2453.2104
$ws.Range("K132").Value = 7359.6312
$ws.Range("M132").Value = -4829.6312
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 2151.5
$ws.Range("J86").Value = 2151.5
$ws.Range("L86").Value = 6454.5
$ws.Range("N86").Value = -8826.5
$ws.Range("H89").Value = 2151.5
$ws.Range("J89").Value = 2151.5
$ws.Range("L89").Value = 19363.5
$ws.Range("N89").Value = -31219.5
$ws.Range("H134").Value = 5570151
$ws.Range("J134").Value = 33216.566
$ws.Range("L134").Value = 99649.698
$ws.Range("N134").Value = -109789.698
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 31242.25
$ws.Range("I57").Value = 0
$ws.Range("J57").Value = 31242.25
$ws.Range("K57").Value = 0
$ws.Range("L57").Value = 31242.25
$ws.Range("M57").Value = $null
$ws.Range("N57").Value = -32882.25
$ws.Range("H102").Value = 3430.8333
$ws.Range("I102").Value = 3430.8333
$ws.Range("K102").Value = 3430.8333
$ws.Range("M102").Value = -1808.8333
$ws.Range("H122").Value = 2295
$ws.Range("I122").Value = 2177.5
$ws.Range("K122").Value = 6532.5
$ws.Range("M122").Value = -4082.5
$ws.Range("H136").Value = 74178.71000000001
$ws.Range("J136").Value = 74178.71000000001
$ws.Range("L136").Value = 222536.13
$ws.Range("N136").Value = -227636.13
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 4862524.5
$ws.Range("I68").Value = 6077509
$ws.Range("K68").Value = 6077509
$ws.Range("M68").Value = -6076760
$ws.Range("H71").Value = 4862524.5
$ws.Range("I71").Value = 6077509
$ws.Range("K71").Value = 30387545
$ws.Range("M71").Value = -30383801
$ws.Range("H93").Value = 1919626.8
$ws.Range("J93").Value = 4279652
$ws.Range("L93").Value = 4279652
$ws.Range("N93").Value = -4282148
$ws.Range("H122").Value = 3579.0293
$ws.Range("I122").Value = 3364.1292
$ws.Range("K122").Value = 10092.3876
$ws.Range("M122").Value = -7642.3876
$ws.Range("H132").Value = 3919
$ws.Range("I132").Value = 2353.125
$ws.Range("J132").Value = 6006.8335
$ws.Range("K132").Value = 7059.375
$ws.Range("L132").Value = 18020.5005
$ws.Range("M132").Value = -4529.375
$ws.Range("N132").Value = -23080.5005
$ws.Range("H136").Value = 3905.6667
$ws.Range("I136").Value = 2686.9
$ws.Range("K136").Value = 8060.700000000001
$ws.Range("M136").Value = -5510.700000000001
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 530024.3
$ws.Range("I100").Value = 1514.1111
$ws.Range("J100").Value = 3701085.8
$ws.Range("K100").Value = 3028.2222
$ws.Range("L100").Value = 7402171.6
$ws.Range("M100").Value = -2487.2222
$ws.Range("N100").Value = -7403253.6
$ws.Range("H122").Value = 3332.5881
$ws.Range("I122").Value = 3034.25
$ws.Range("K122").Value = 9102.75
$ws.Range("M122").Value = -6652.75
$ws.Range("H126").Value = 8543.333000000001
$ws.Range("I126").Value = 8325.385
$ws.Range("J126").Value = 9960
$ws.Range("K126").Value = 24976.155
$ws.Range("L126").Value = 29880
$ws.Range("M126").Value = -22506.155
$ws.Range("N126").Value = -34820
$ws.Range("H132").Value = 1252158
$ws.Range("I132").Value = 2377.3333
$ws.Range("K132").Value = 7131.999899999999
$ws.Range("M132").Value = -4601.999899999999
$ws.Range("H136").Value = 537574.25
$ws.Range("I136").Value = 12347.706
$ws.Range("K136").Value = 37043.118
$ws.Range("M136").Value = -34493.118
